$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Price" (D) and "Volume(1h)" (E) columns with freshly scraped values.
# Some Price values (e.g. "312.70", "5.000") are strings that look like plain numbers;
# force the cell to text first so Excel keeps the exact original formatting/precision,
# then restore the default "Normal" style so no stray number formatting is left behind.

$ws.Range("D2").Value = '27.514.53'
$ws.Range("E2").Value = '  +2.20%  '
$ws.Range("D3").Value = '1.869.59'
$ws.Range("E3").Value = '  +1.31%  '
$ws.Range("E4").Value = '  +0.61%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.06%  '
$ws.Range("E6").Value = '  +0.62%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4779'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.44%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3775'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07373'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9369'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.75'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.14%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07846'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.99%  '
$ws.Range("D13").Value = '1.871.14'
$ws.Range("E13").Value = '  -0.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.439'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.582'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.76%  '
$ws.Range("E16").Value = '  +2.29%  '
$ws.Range("E17").Value = '  +0.57%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008912'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.013'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.95'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.68%  '
$ws.Range("D21").Value = '27.517.32'
$ws.Range("E21").Value = '  +2.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.138'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.74'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.964'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.05'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.08%  '
$ws.Range("E26").Value = '  +2.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.023'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '116.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.000'
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = '  +0.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.338'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.220'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7550'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.89%  '
$ws.Range("E34").Value = '  +2.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.692'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02050'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.79%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.118'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05278'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.999'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.04%  '
$ws.Range("E40").Value = '  +2.89%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.075'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.33%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1525'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.453'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.64'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4819'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.014'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.660'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.46%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.74'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '67.53'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.19%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06092'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9286'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.91%  '
